$d = $word.ActiveDocument
$dash = [char]0x2013

# --- 1. Extend the first paragraph: two trailing spaces, then a red parenthetical note ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertAfter("  ")

# -- run: "(This is a change – Ve" (red) --
$text1 = "(This is a change " + $dash + " Ve"
$pos = $d.Paragraphs(1).Range.End - 1
$d.Paragraphs(1).Range.InsertAfter($text1)
$r1 = $d.Range($pos, $pos + $text1.Length)
$r1.Font.Color = 255

# -- run: "rsion for main branch" (red) --
$text2 = "rsion for main branch"
$pos = $d.Paragraphs(1).Range.End - 1
$d.Paragraphs(1).Range.InsertAfter($text2)
$r2 = $d.Range($pos, $pos + $text2.Length)
$r2.Font.Color = 255

# -- run: ")" (red) --
$text3 = ")"
$pos = $d.Paragraphs(1).Range.End - 1
$d.Paragraphs(1).Range.InsertAfter($text3)
$r3 = $d.Range($pos, $pos + $text3.Length)
$r3.Font.Color = 255

Write-Host "paragraph 1 updated"
